$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old "Prediksi Suhu" column to C)
$ws.Columns.Item(2).Insert()

# New header for inserted column B, and keep header style consistent (same style as A1/old B1)
$ws.Range("B1").Value = "Nh Label"
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats

# New "Nh Label" integer values for B2:B10
$nhValues = @(2, 2, 2, 3, 0, 8, 2, 2, 2)
for ($i = 0; $i -lt $nhValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $nhValues[$i]
}

# Updated "Prediksi Suhu" values now living in column C
$suhuValues = @(29.7065722, 29.7065722, 29.7065722, 29.70656551, 29.70658557, 29.70653207, 29.7065722, 29.7065722, 29.7065722)
for ($i = 0; $i -lt $suhuValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $suhuValues[$i]
}
